$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1349983333333333
$ws.Range("H2").Value = 0.404995
$ws.Range("I2").Value = 0.06188478316908706
$ws.Range("J2").Value = 0.06188478316908706
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.155697
$ws.Range("N2").Value = 0.467091
$ws.Range("O2").Value = 0.021288392311201
$ws.Range("P2").Value = 0.021288392311201
$ws.Range("Q2").Value = 0.021018835505
$ws.Range("R2").Value = 0.189169519545
$ws.Range("S2").Value = 0.001317427542197134
$ws.Range("T2").Value = 0.001317427542197134

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1349983333333333
$ws.Range("H3").Value = 0.404995
$ws.Range("I3").Value = 0.06188478316908706
$ws.Range("J3").Value = 0.06188478316908706
$ws.Range("O3").Value = 0.01134295290047287
$ws.Range("P3").Value = 0.01134295290047287
$ws.Range("Q3").Value = 0.011199326735
$ws.Range("R3").Value = 0.100793940615
$ws.Range("S3").Value = 0.0007019561807429304
$ws.Range("T3").Value = 0.0007019561807429304

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1349983333333333
$ws.Range("H4").Value = 0.404995
$ws.Range("I4").Value = 0.06188478316908706
$ws.Range("J4").Value = 0.06188478316908706
$ws.Range("M4").Value = 7.075048
$ws.Range("N4").Value = 21.225144
$ws.Range("O4").Value = 0.9673686547883261
$ws.Range("P4").Value = 0.9673686547883261
$ws.Range("Q4").Value = 0.9551196882533333
$ws.Range("R4").Value = 8.59607719428
$ws.Range("S4").Value = 0.05986539944614699
$ws.Range("T4").Value = 0.05986539944614699

$ws.Range("I5").Value = 0.4284959871424753
$ws.Range("J5").Value = 0.4284959871424753
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.155697
$ws.Range("N5").Value = 0.467091
$ws.Range("O5").Value = 0.021288392311201
$ws.Range("P5").Value = 0.021288392311201
$ws.Range("Q5").Value = 0.145536369477
$ws.Range("R5").Value = 1.309827325293
$ws.Range("S5").Value = 0.009121990678064355
$ws.Range("T5").Value = 0.009121990678064355

$ws.Range("I6").Value = 0.4284959871424753
$ws.Range("J6").Value = 0.4284959871424753
$ws.Range("O6").Value = 0.01134295290047287
$ws.Range("P6").Value = 0.01134295290047287
$ws.Range("S6").Value = 0.004860409800198725
$ws.Range("T6").Value = 0.004860409800198725

$ws.Range("I7").Value = 0.4284959871424753
$ws.Range("J7").Value = 0.4284959871424753
$ws.Range("M7").Value = 7.075048
$ws.Range("N7").Value = 21.225144
$ws.Range("O7").Value = 0.9673686547883261
$ws.Range("P7").Value = 0.9673686547883261
$ws.Range("Q7").Value = 6.613337442568
$ws.Range("R7").Value = 59.520036983112
$ws.Range("S7").Value = 0.4145135866642122
$ws.Range("T7").Value = 0.4145135866642122

$ws.Range("G8").Value = 1.111707
$ws.Range("H8").Value = 3.335121
$ws.Range("I8").Value = 0.5096192296884376
$ws.Range("J8").Value = 0.5096192296884376
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.155697
$ws.Range("N8").Value = 0.467091
$ws.Range("O8").Value = 0.021288392311201
$ws.Range("P8").Value = 0.021288392311201
$ws.Range("Q8").Value = 0.173089444779
$ws.Range("R8").Value = 1.557805003011
$ws.Range("S8").Value = 0.01084897409093951
$ws.Range("T8").Value = 0.01084897409093951

$ws.Range("G9").Value = 1.111707
$ws.Range("H9").Value = 3.335121
$ws.Range("I9").Value = 0.5096192296884376
$ws.Range("J9").Value = 0.5096192296884376
$ws.Range("O9").Value = 0.01134295290047287
$ws.Range("P9").Value = 0.01134295290047287
$ws.Range("Q9").Value = 0.09222610101300001
$ws.Range("R9").Value = 0.8300349091170001
$ws.Range("S9").Value = 0.005780586919531211
$ws.Range("T9").Value = 0.005780586919531211

$ws.Range("G10").Value = 1.111707
$ws.Range("H10").Value = 3.335121
$ws.Range("I10").Value = 0.5096192296884376
$ws.Range("J10").Value = 0.5096192296884376
$ws.Range("M10").Value = 7.075048
$ws.Range("N10").Value = 21.225144
$ws.Range("O10").Value = 0.9673686547883261
$ws.Range("P10").Value = 0.9673686547883261
$ws.Range("Q10").Value = 7.865380386936
$ws.Range("R10").Value = 70.788423482424
$ws.Range("S10").Value = 0.4929896686779669
$ws.Range("T10").Value = 0.4929896686779669
